# DPLKINV125-001 - Profile Reksa - Investasi - Reksadana Hapus Data.xlsx
# Regression pass: bump the generated "Kode Reksadana" sample value from
# RD00014 to RD00015 (used both in the PREPARATION instructions cell and
# in the dedicated KODE_REKSADANA column), and move the active selection
# over to the KODE_REKSADANA cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 (PREPARATION): update the "Kode Reksadana : ..." line at the end of
# the multi-line preparation text.
$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 18/19/20/21 - Pimpinan Kelompok Investasi/Pengelolan Investasi/Analis;`nKode Reksadana : RD00015"

# N2 (KODE_REKSADANA): the generated code itself.
$ws.Range("N2").Value = "RD00015"

# Move/scroll the view so column O (KODE_REKSADANA) is visible and select
# O2 as the active cell, matching the tester's new focus point.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("O2").Select()
